$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Existing-cell value tweaks (rows 11, 13, 14, 15, 17 - column B)
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = 2
$ws.Range("B13").Value = 2
$ws.Range("B14").Value = 2
$ws.Range("B15").Value = 0
$ws.Range("B17").Value = 2

# ---------------------------------------------------------------------------
# 2) Highlight fills on existing cells.
#    Order matters: it determines the order new <fill>/<xf> records are
#    appended to styles.xml, so do the "yellow over a white-filled cell"
#    highlight first (creates fgColor=yellow/bgColor=black), then the new
#    row's red highlight (creates fgColor=red/bgColor=black), then the
#    "yellow over a no-fill cell" highlight (creates fgColor=yellow/bgColor
#    automatic) last.
# ---------------------------------------------------------------------------
$ws.Range("D15").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 3) New row 23 - day-of-month header strip (values 0..40), with a subset of
#    cells (roughly the Sundays in the calendar) highlighted red.
# ---------------------------------------------------------------------------
$row23 = [ordered]@{
    "B23"  = 0
    "C23"  = 1
    "D23"  = 2
    "E23"  = 3
    "F23"  = 4
    "G23"  = 5
    "H23"  = 6
    "I23"  = 7
    "J23"  = 8
    "K23"  = 9
    "L23"  = 10
    "M23"  = 11
    "N23"  = 12
    "O23"  = 13
    "P23"  = 14
    "Q23"  = 15
    "R23"  = 16
    "S23"  = 17
    "T23"  = 18
    "U23"  = 19
    "V23"  = 20
    "W23"  = 21
    "X23"  = 22
    "Y23"  = 23
    "Z23"  = 24
    "AA23" = 25
    "AB23" = 26
    "AC23" = 27
    "AD23" = 28
    "AE23" = 29
    "AF23" = 30
    "AG23" = 31
    "AH23" = 32
    "AI23" = 33
    "AJ23" = 34
    "AK23" = 35
    "AL23" = 36
    "AM23" = 37
    "AN23" = 38
    "AO23" = 39
    "AP23" = 40
}
foreach ($addr in $row23.Keys) {
    $ws.Range($addr).Value = $row23[$addr]
}

$row23Red = @("B23","L23","V23","Y23","Z23","AB23","AC23","AE23","AH23","AI23","AL23","AO23")
foreach ($addr in $row23Red) {
    $ws.Range($addr).Interior.Color = 255
}

# ---------------------------------------------------------------------------
# 4) Remaining yellow highlights - reuse the "yellow / automatic background"
#    fill created below for the cells that had no fill beforehand, and the
#    "yellow / black background" fill for the cells that already had a white
#    fill beforehand.
# ---------------------------------------------------------------------------
$ws.Range("C11").Interior.Color = 65535
$ws.Range("C12").Interior.Color = 65535
$ws.Range("C13").Interior.Color = 65535
$ws.Range("C14").Interior.Color = 65535

$ws.Range("D16").Interior.Color = 65535
$ws.Range("D17").Interior.Color = 65535

# E15 loses its old fill (becomes the plain "no fill" border style instead).
$ws.Range("E15").Interior.ColorIndex = -4142

# ---------------------------------------------------------------------------
# 5) Selection moves to X23 (the new row).
# ---------------------------------------------------------------------------
$null = $ws.Range("X23").Select()
